$d = $word.ActiveDocument

# --- 1. Un-bold the tail of the final paragraph (paragraph mark + everything from the space after
#        "...completion rate target." through to the end of the paragraph). Do this by clearing
#        bold across the whole paragraph (which also clears the paragraph-mark's own bold in pPr/rPr)
#        and then re-applying bold to the portion that should stay bold. Do this FIRST, while the
#        paragraph's run layout is still in its original (finely split) state. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Font.Bold = $false

$keepBold = $d.Content
$keepBold.Find.Execute(" target.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$reboldRange = $d.Range($lastPara.Range.Start, $keepBold.End)
$reboldRange.Font.Bold = $true

# --- 2. Update creation/last-updated dates: 12/18/2018 -> 12/17/2019 ---
$d.Content.Find.Execute("on 12/18/2018. It was last updated on 12/18/2018.", $true, $false, $false, $false, $false, $true, 1, $false, "on 12/17/2019. It was last updated on 12/17/2019.", 2)

# --- 3. "goals during 2018" -> "goals during 2019" ---
$d.Content.Find.Execute("goals during 2018", $true, $false, $false, $false, $false, $true, 1, $false, "goals during 2019", 2)

# --- 4. "goals that were met for 2018" -> "goals that were met for 2019" ---
$d.Content.Find.Execute("goals that were met for 2018", $true, $false, $false, $false, $false, $true, 1, $false, "goals that were met for 2019", 2)

# --- 5. Merge the two "70% of youth will report stable..." runs into one (no text change, just
#        collapses the trailing-space run into the sentence run) ---
$d.Content.Find.Execute("70% of youth will report stable or improved attitudes toward marijuana use at post-test ", $true, $false, $false, $false, $false, $true, 1, $false, "70% of youth will report stable or improved attitudes toward marijuana use at post-test ", 2)

# --- 6. "In 2018, Campus Connections met" -> "In 2019, Campus Connections met" ---
$d.Content.Find.Execute("In 2018, Campus Connections met", $true, $false, $false, $false, $false, $true, 1, $false, "In 2019, Campus Connections met", 2)

# --- 7. "A total of 226 youth participants" -> "A total of 228 youth participants" ---
$d.Content.Find.Execute("A total of 226 youth participants", $true, $false, $false, $false, $false, $true, 1, $false, "A total of 228 youth participants", 2)

# --- 8. "...was 92% and exceeded" -> "...was 91% and exceeded" (only touch this occurrence, not "for the 2018 year") ---
$d.Content.Find.Execute("for the 2018 year was 92%", $true, $false, $false, $false, $false, $true, 1, $false, "for the 2018 year was 91%", 2)

# --- 9. Move the "_GoBack" bookmark from the end of the document to right after " goal was met." ---
$metRange = $d.Content
$metRange.Find.Execute("goal was met.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$metRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $metRange)

Write-Output "done"
